# Update "想去人数" (want-to-go count) figures and one sold-out flag
# across the four sheets: 展览(1), 演出(2), 本地生活(3), 全部类型(4).
#
# Values below are the NEW values for column F (and the one special
# column G change on 演出 row 15), keyed by row number.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)   # 展览
$sheet2 = $wb.Worksheets.Item(2)   # 演出
$sheet3 = $wb.Worksheets.Item(3)   # 本地生活
$sheet4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 -----------------------------------------------------
$sheet1.Range("F2").Value  = 1526
$sheet1.Range("F5").Value  = 7880
$sheet1.Range("F6").Value  = 4902
$sheet1.Range("F7").Value  = 7194
$sheet1.Range("F8").Value  = 298
$sheet1.Range("F9").Value  = 1531
$sheet1.Range("F12").Value = 73
$sheet1.Range("F13").Value = 1189
$sheet1.Range("F14").Value = 194
$sheet1.Range("F15").Value = 566
$sheet1.Range("F16").Value = 30
$sheet1.Range("F20").Value = 1245
$sheet1.Range("F21").Value = 964
$sheet1.Range("F23").Value = 48
$sheet1.Range("F24").Value = 1269
$sheet1.Range("F25").Value = 54
$sheet1.Range("F26").Value = 166
$sheet1.Range("F28").Value = 20
$sheet1.Range("F30").Value = 226
$sheet1.Range("F31").Value = 1025
$sheet1.Range("F33").Value = 14
$sheet1.Range("F35").Value = 137
$sheet1.Range("F36").Value = 42
$sheet1.Range("F37").Value = 559
$sheet1.Range("F38").Value = 576
$sheet1.Range("F39").Value = 425
$sheet1.Range("F40").Value = 91
$sheet1.Range("F41").Value = 64
$sheet1.Range("F42").Value = 110
$sheet1.Range("F43").Value = 438
$sheet1.Range("F44").Value = 1207
$sheet1.Range("F45").Value = 611
$sheet1.Range("F46").Value = 169

# --- Sheet 2: 演出 -------------------------------------------------------
$sheet2.Range("F5").Value  = 20

# G15 switches from a numeric price (380) to a text flag meaning
# "not sellable" - must be written as text, not a number.
$sheet2.Range("G15").Value = "不可售"

$sheet2.Range("F28").Value = 2
$sheet2.Range("F29").Value = 38
$sheet2.Range("F32").Value = 879
$sheet2.Range("F33").Value = 7
$sheet2.Range("F34").Value = 1004
$sheet2.Range("F41").Value = 109
$sheet2.Range("F43").Value = 19
$sheet2.Range("F46").Value = 84
$sheet2.Range("F48").Value = 10

# --- Sheet 3: 本地生活 ---------------------------------------------------
$sheet3.Range("F5").Value  = 865
$sheet3.Range("F7").Value  = 198
$sheet3.Range("F8").Value  = 108
$sheet3.Range("F9").Value  = 1768
$sheet3.Range("F10").Value = 2672

# --- Sheet 4: 全部类型 ---------------------------------------------------
$sheet4.Range("F3").Value  = 1526
$sheet4.Range("F4").Value  = 865
$sheet4.Range("F7").Value  = 7880
$sheet4.Range("F8").Value  = 198
$sheet4.Range("F9").Value  = 4902
$sheet4.Range("F10").Value = 7194
$sheet4.Range("F11").Value = 298
$sheet4.Range("F12").Value = 1531
$sheet4.Range("F14").Value = 108
$sheet4.Range("F16").Value = 1768
$sheet4.Range("F17").Value = 2672
$sheet4.Range("F19").Value = 1189
$sheet4.Range("F20").Value = 194
$sheet4.Range("F22").Value = 30
$sheet4.Range("F24").Value = 1245
$sheet4.Range("F26").Value = 964
$sheet4.Range("F27").Value = 1269
$sheet4.Range("F28").Value = 166
$sheet4.Range("F29").Value = 20
$sheet4.Range("F30").Value = 226
$sheet4.Range("F33").Value = 38
$sheet4.Range("F34").Value = 879
$sheet4.Range("F36").Value = 7
$sheet4.Range("F37").Value = 137
$sheet4.Range("F38").Value = 1004
$sheet4.Range("F39").Value = 576
$sheet4.Range("F41").Value = 91
$sheet4.Range("F42").Value = 110
$sheet4.Range("F44").Value = 438
$sheet4.Range("F45").Value = 611
$sheet4.Range("F46").Value = 109
